# Calculator.xlsx — "Added PSRAM numbers to calc"
#
# This script updates the live LCD calculator (rows 2-10) to use the
# "800x480" LCD preset as its active input (instead of the 1024x600 one),
# moves the old 1024x600 numbers into a third preset column (J:L, rows
# 25-33), and adds a brand-new "PSRAM Specs" block (columns V:X) with
# derived PSRAM transfer-rate / size figures. It also appends a note
# about PSRAM transfer-rate overhead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Preserve the old "live" LCD numbers (1024x600...) as a new preset
#    in columns J:L (rows 25-33), mirroring the B:D / F:H preset blocks.
# ---------------------------------------------------------------------
$ws.Range("J25:L25").Merge()
$ws.Range("J25").Value = "LCD Specs"

$ws.Range("J26").Value = "Width"
$ws.Range("K26").Value = 1024
$ws.Range("L26").Value = "Px"

$ws.Range("J27").Value = "Height"
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = "Px"

$ws.Range("J28").Value = "Refresh"
$ws.Range("K28").Value = 60
$ws.Range("L28").Value = "Hz"

$ws.Range("J29").Value = "Total Width"
$ws.Range("K29").Value = 1114
$ws.Range("L29").Value = "Px"

$ws.Range("J30").Value = "Total Height"
$ws.Range("K30").Value = 610
$ws.Range("L30").Value = "Px"

$ws.Range("J31").Value = "BPP"
$ws.Range("K31").Value = 18
$ws.Range("L31").Value = "bits"

$ws.Range("J32").Value = "Hsync"
$ws.Range("K32").Formula = "=K29-K26"
$ws.Range("L32").Value = "Px"

$ws.Range("J33").Value = "Vsync"
$ws.Range("K33").Formula = "=K30-K27"
$ws.Range("L33").Value = "Px"

# Copy header styling/format from the existing preset header cells.
$ws.Range("B25").Copy()
$ws.Range("J25").PasteSpecial(-4122)
$ws.Range("J25:L25").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 2. Switch the live calculator's LCD inputs to the 800x480 preset.
# ---------------------------------------------------------------------
$ws.Range("C3").Value = 800
$ws.Range("C4").Value = 480
$ws.Range("C6").Value = 816
$ws.Range("C7").Value = 496

# ---------------------------------------------------------------------
# 3. Relabel the two "RAM->HSTX" hardware-requirement rates that are now
#    generic to any memory (PSRAM included), and add the new PSRAM
#    Specs block.
# ---------------------------------------------------------------------
$ws.Range("N7").Value = "Mem->HSTX Rate(Full)"
$ws.Range("N8").Value = "Mem->HSTX Rate (565)"

$ws.Range("V2:X2").Merge()
$ws.Range("V2").Value = "PSRAM Specs"
$ws.Range("V2:X2").HorizontalAlignment = -4108

$ws.Range("V3").Value = "Max Transfer Rate"
$ws.Range("W3").Formula = "=4*S3"
$ws.Range("X3").Value = "Mbps"

$ws.Range("V4").Value = "Max Size"
$ws.Range("W4").Value = 128
$ws.Range("X4").Value = "Mbit"

$ws.Range("V5").Value = "Max Size"
$ws.Range("W5").Formula = "=(W4/8)*1024"
$ws.Range("X5").Value = "Kb"

# ---------------------------------------------------------------------
# 4. Append the PSRAM caveat to the notes box.
# ---------------------------------------------------------------------
$ws.Range("B18").Value = "Notes:`n-DDR only works with 16 BPP (i.e. 565) as it uses the same 32 Bit buffer for 2 px`n-Actual PSRAM transfer rate will be slightly slower due to overhead, negligible for large read/writes, but high for random r/w"

# ---------------------------------------------------------------------
# 5. Cosmetic: column widths / view matching the refreshed layout.
# ---------------------------------------------------------------------
$ws.Range("N1").ColumnWidth = 19.55
$ws.Range("V1").ColumnWidth = 15.33
$ws.Range("W1").ColumnWidth = 8.89
$ws.Range("X1").ColumnWidth = 7.33

$ws.Range("G6").Select()
$excel.ActiveWindow.Zoom = 100
